$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") -- copy H1's format
# (bold/border/center/top-align) onto I1:J1 so they reuse the existing
# header style, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..21 for columns I (col 9) and J (col 10)
$values = @(
    @(3,3),
    @(6,6),
    @(6,7),
    @(9,9),
    @(6,6),
    @(8,8),
    @(6,6),
    @(6,6),
    @(5,6),
    @(6,7),
    @(8,8),
    @(5,6),
    @(8,8),
    @(9,9),
    @(2,2),
    @(7,7),
    @(7,7),
    @(8,8),
    @(6,6),
    @(3,3)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
